$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")
$ws.Range("H1").Value = "TestValue"
Write-Output ($ws.Range("H1").Text)
